$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row: date serial 46001 (12/10/2025) with an error count
# of 0. Copy the format of the row above first so the new date cell picks
# up the same date number format used by the rest of column A, then set
# the actual values.
$ws.Range("A31").Copy()
$null = $ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A32").Value = 46001
$ws.Range("B32").Value = 0

# Leave the selection on the next blank cell below the data that was
# just entered.
$null = $ws.Range("B33").Select()
